$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the value columns as Text so the numeric-looking strings are stored
# as shared strings (matching the source data), not auto-converted to numbers.
$ws.Range("D2:F16").NumberFormat = "@"
$ws.Range("J2:M16").NumberFormat = "@"

$ws.Range("D2").Value = "-10.9932194001311"
$ws.Range("E2").Value = "64.0370307305784"
$ws.Range("F2").Value = "2.19856237737517e-16"
$ws.Range("J2").Value = "-34.5456291467162"
$ws.Range("K2").Value = "-23.9209739212368"
$ws.Range("L2").Value = "17.5"
$ws.Range("M2").Value = "46.7333015339765"

$ws.Range("D3").Value = "-11.6321953539881"
$ws.Range("E3").Value = "64.5925315019376"
$ws.Range("F3").Value = "1.70327372671198e-17"
$ws.Range("J3").Value = "-35.6551707539455"
$ws.Range("K3").Value = "-25.2048339055072"
$ws.Range("L3").Value = "17.5"
$ws.Range("M3").Value = "47.9300023297263"

$ws.Range("D4").Value = "-15.0456482614919"
$ws.Range("E4").Value = "65.3128771122382"
$ws.Range("F4").Value = "6.63883335973362e-23"
$ws.Range("J4").Value = "-39.2004817069638"
$ws.Range("K4").Value = "-30.0138978805907"
$ws.Range("L4").Value = "17.5"
$ws.Range("M4").Value = "52.1071897937773"

$ws.Range("D5").Value = "-14.0743090735612"
$ws.Range("E5").Value = "64.0959410014433"
$ws.Range("F5").Value = "2.81537071135295e-21"
$ws.Range("J5").Value = "-35.8507432583315"
$ws.Range("K5").Value = "-26.9385839261492"
$ws.Range("L5").Value = "17.5"
$ws.Range("M5").Value = "48.8946635922404"

$ws.Range("D6").Value = "8.67946751904764"
$ws.Range("E6").Value = "34.0561802025822"
$ws.Range("F6").Value = "3.79206722649747e-10"
$ws.Range("J6").Value = "11.4430034686036"
$ws.Range("K6").Value = "18.4393494725729"
$ws.Range("L6").Value = "17.5"
$ws.Range("M6").Value = "2.55882352941176"

$ws.Range("D7").Value = "-0.420955880210163"
$ws.Range("E7").Value = "65.9475693641348"
$ws.Range("F7").Value = "0.675156044280717"
$ws.Range("J7").Value = "-6.87265354461256"
$ws.Range("K7").Value = "4.47925195311279"
$ws.Range("L7").Value = "46.7333015339765"
$ws.Range("M7").Value = "47.9300023297263"

$ws.Range("D8").Value = "-2.10316666421769"
$ws.Range("E8").Value = "61.4299622560178"
$ws.Range("F8").Value = "0.0395545430011091"
$ws.Range("J8").Value = "-10.482486160765"
$ws.Range("K8").Value = "-0.265290358836554"
$ws.Range("L8").Value = "46.7333015339765"
$ws.Range("M8").Value = "52.1071897937773"

$ws.Range("D9").Value = "-0.867058542091472"
$ws.Range("E9").Value = "59.2595615280683"
$ws.Range("F9").Value = "0.38940734031658"
$ws.Range("J9").Value = "-7.14889053208308"
$ws.Range("K9").Value = "2.82616641555524"
$ws.Range("L9").Value = "46.7333015339765"
$ws.Range("M9").Value = "48.8946635922404"

$ws.Range("D10").Value = "21.5512730372604"
$ws.Range("E10").Value = "33.7415444319119"
$ws.Range("F10").Value = "2.64877596791972e-21"
$ws.Range("J10").Value = "40.0077322132428"
$ws.Range("K10").Value = "48.3412237958865"
$ws.Range("L10").Value = "46.7333015339765"
$ws.Range("M10").Value = "2.55882352941176"

$ws.Range("D11").Value = "-1.66411568527365"
$ws.Range("E11").Value = "62.22079215284"
$ws.Range("F11").Value = "0.101118166962509"
$ws.Range("J11").Value = "-9.19456063486408"
$ws.Range("K11").Value = "0.840185706762245"
$ws.Range("L11").Value = "47.9300023297263"
$ws.Range("M11").Value = "52.1071897937773"

$ws.Range("D12").Value = "-0.394283509893122"
$ws.Range("E12").Value = "60.1441661965777"
$ws.Range("F12").Value = "0.694766337294604"
$ws.Range("J12").Value = "-5.85838531125264"
$ws.Range("K12").Value = "3.92906278622457"
$ws.Range("L12").Value = "47.9300023297263"
$ws.Range("M12").Value = "48.8946635922404"

$ws.Range("D13").Value = "22.7609488720519"
$ws.Range("E13").Value = "33.7845636403643"
$ws.Range("F13").Value = "4.52153999486825e-22"
$ws.Range("J13").Value = "41.3191952122571"
$ws.Range("K13").Value = "49.423162388372"
$ws.Range("L13").Value = "47.9300023297263"
$ws.Range("M13").Value = "2.55882352941176"

$ws.Range("D14").Value = "1.52578747925203"
$ws.Range("E14").Value = "65.6685736360082"
$ws.Range("F14").Value = "0.131864961649821"
$ws.Range("J14").Value = "-0.991610920912855"
$ws.Range("K14").Value = "7.41666332398662"
$ws.Range("L14").Value = "52.1071897937773"
$ws.Range("M14").Value = "48.8946635922404"

$ws.Range("D15").Value = "31.8463419756661"
$ws.Range("E15").Value = "34.2974374030118"
$ws.Range("F15").Value = "4.64458789074224e-27"
$ws.Range("J15").Value = "46.38749335084"
$ws.Range("K15").Value = "52.709239177891"
$ws.Range("L15").Value = "52.1071897937773"
$ws.Range("M15").Value = "2.55882352941176"

$ws.Range("D16").Value = "31.9310224417696"
$ws.Range("E16").Value = "34.4956896545343"
$ws.Range("F16").Value = "3.32481233563374e-27"
$ws.Range("J16").Value = "43.388363485552"
$ws.Range("K16").Value = "49.2833166401053"
$ws.Range("L16").Value = "48.8946635922404"
$ws.Range("M16").Value = "2.55882352941176"
